# TC_144_145 - Updated test data for 5,24,40V,BatteryStandby and AC Calculations test cases
#
# The "Add Panels" sheet holds a single CPU/PSU change scenario (row 8).
# This edit:
#   - adds a reference note in B4 ("NGC-488/T386 OR TC-144 & 145")
#   - renames the panel under test from "MX1000" to "MX 1000"
#   - blanks out the now-unused CPU-type columns (C8, F8)
#   - renames the PSU part number from "PSU830" to "PSB800"
#   - replaces the old numeric Alarm/Standby current readings with the
#     new (text) readings: 0.250 / 0.562 / 0.235 / 0.547
#   - leaves the selection on B8, scrolled back to column A

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# Row 4 - add the reference/user-story note; the cell had no content or
# fill before, so drop its border formatting too (plain/default style).
$ws.Range("B4").ClearFormats()
$ws.Range("B4").Value = "NGC-488/T386 OR TC-144 & 145"

# Row 8 - panel/CPU/PSU identifiers
$ws.Range("A8").Value = "MX 1000"
$ws.Range("C8").Value = "'"
$ws.Range("F8").Value = "'"
$ws.Range("G8").Value = "PSB800"

# Row 8 - Alarm/Standby current readings (now stored as text values)
$ws.Range("H8").Value = "'0.250"
$ws.Range("I8").Value = "'0.562"
$ws.Range("J8").Value = "'0.235"
$ws.Range("K8").Value = "'0.547"

# Restore the view: select B8 and scroll the sheet back to show column A.
[void]$ws.Range("A1").Select()
[void]$ws.Range("B8").Select()
